$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "1" -as [double]
$ws.Range("F2").Value = "0.3333333333333333" -as [double]
$ws.Range("G2").Value = "0.01675466666666667" -as [double]
$ws.Range("H2").Value = "0.050264" -as [double]
$ws.Range("I2").Value = "0.0001854906931657378" -as [double]
$ws.Range("J2").Value = "0.0001854906931657378" -as [double]
$ws.Range("M2").Value = "0.4102596666666667" -as [double]
$ws.Range("N2").Value = "1.230779" -as [double]
$ws.Range("O2").Value = "0.003499619873322347" -as [double]
$ws.Range("P2").Value = "0.003499619873322347" -as [double]
$ws.Range("Q2").Value = "0.006873763961777779" -as [double]
$ws.Range("R2").Value = "0.06186387565600001" -as [double]
$ws.Range("S2").Value = "6.491469161191537E-07" -as [double]
$ws.Range("T2").Value = "6.491469161191535E-07" -as [double]
$ws.Range("E3").Value = "1" -as [double]
$ws.Range("F3").Value = "0.3333333333333333" -as [double]
$ws.Range("G3").Value = "0.01675466666666667" -as [double]
$ws.Range("H3").Value = "0.050264" -as [double]
$ws.Range("I3").Value = "0.0001854906931657378" -as [double]
$ws.Range("J3").Value = "0.0001854906931657378" -as [double]
$ws.Range("O3").Value = "0.8692174743460166" -as [double]
$ws.Range("P3").Value = "0.8692174743460165" -as [double]
$ws.Range("Q3").Value = "1.707269922557333" -as [double]
$ws.Range("R3").Value = "15.365429303016" -as [double]
$ws.Range("S3").Value = "0.0001612317518282145" -as [double]
$ws.Range("T3").Value = "0.0001612317518282145" -as [double]
$ws.Range("E4").Value = "1" -as [double]
$ws.Range("F4").Value = "0.3333333333333333" -as [double]
$ws.Range("G4").Value = "0.01675466666666667" -as [double]
$ws.Range("H4").Value = "0.050264" -as [double]
$ws.Range("I4").Value = "0.0001854906931657378" -as [double]
$ws.Range("J4").Value = "0.0001854906931657378" -as [double]
$ws.Range("N4").Value = "44.764041" -as [double]
$ws.Range("O4").Value = "0.1272829057806611" -as [double]
$ws.Range("P4").Value = "0.1272829057806611" -as [double]
$ws.Range("Q4").Value = "0.2500021952026666" -as [double]
$ws.Range("R4").Value = "2.250019756824" -as [double]
$ws.Range("S4").Value = "2.360979442140413E-05" -as [double]
$ws.Range("T4").Value = "2.360979442140413E-05" -as [double]
$ws.Range("I5").Value = "0.9933938536206305" -as [double]
$ws.Range("J5").Value = "0.9933938536206304" -as [double]
$ws.Range("M5").Value = "0.4102596666666667" -as [double]
$ws.Range("N5").Value = "1.230779" -as [double]
$ws.Range("O5").Value = "0.003499619873322347" -as [double]
$ws.Range("P5").Value = "0.003499619873322347" -as [double]
$ws.Range("Q5").Value = "36.81238532419433" -as [double]
$ws.Range("R5").Value = "331.311467917749" -as [double]
$ws.Range("S5").Value = "0.00347650087216703" -as [double]
$ws.Range("T5").Value = "0.003476500872167029" -as [double]
$ws.Range("I6").Value = "0.9933938536206305" -as [double]
$ws.Range("J6").Value = "0.9933938536206304" -as [double]
$ws.Range("O6").Value = "0.8692174743460166" -as [double]
$ws.Range("P6").Value = "0.8692174743460165" -as [double]
$ws.Range("S6").Value = "0.8634752964749809" -as [double]
$ws.Range("T6").Value = "0.8634752964749808" -as [double]
$ws.Range("I7").Value = "0.9933938536206305" -as [double]
$ws.Range("J7").Value = "0.9933938536206304" -as [double]
$ws.Range("N7").Value = "44.764041" -as [double]
$ws.Range("O7").Value = "0.1272829057806611" -as [double]
$ws.Range("P7").Value = "0.1272829057806611" -as [double]
$ws.Range("S7").Value = "0.1264420562734826" -as [double]
$ws.Range("T7").Value = "0.1264420562734826" -as [double]
$ws.Range("G8").Value = "0.5799533333333334" -as [double]
$ws.Range("I8").Value = "0.006420655686203657" -as [double]
$ws.Range("J8").Value = "0.006420655686203655" -as [double]
$ws.Range("M8").Value = "0.4102596666666667" -as [double]
$ws.Range("N8").Value = "1.230779" -as [double]
$ws.Range("O8").Value = "0.003499619873322347" -as [double]
$ws.Range("P8").Value = "0.003499619873322347" -as [double]
$ws.Range("Q8").Value = "0.2379314612155556" -as [double]
$ws.Range("R8").Value = "2.14138315094" -as [double]
$ws.Range("S8").Value = "2.246985423919845E-05" -as [double]
$ws.Range("T8").Value = "2.246985423919844E-05" -as [double]
$ws.Range("G9").Value = "0.5799533333333334" -as [double]
$ws.Range("I9").Value = "0.006420655686203657" -as [double]
$ws.Range("J9").Value = "0.006420655686203655" -as [double]
$ws.Range("O9").Value = "0.8692174743460166" -as [double]
$ws.Range("P9").Value = "0.8692174743460165" -as [double]
$ws.Range("Q9").Value = "59.09618509192668" -as [double]
$ws.Range("R9").Value = "531.8656658273401" -as [double]
$ws.Range("S9").Value = "0.005580946119207333" -as [double]
$ws.Range("T9").Value = "0.00558094611920733" -as [double]
$ws.Range("G10").Value = "0.5799533333333334" -as [double]
$ws.Range("I10").Value = "0.006420655686203657" -as [double]
$ws.Range("J10").Value = "0.006420655686203655" -as [double]
$ws.Range("N10").Value = "44.764041" -as [double]
$ws.Range("O10").Value = "0.1272829057806611" -as [double]
$ws.Range("P10").Value = "0.1272829057806611" -as [double]
$ws.Range("S10").Value = "0.0008172397127571263" -as [double]
$ws.Range("T10").Value = "0.0008172397127571261" -as [double]
